$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 74, pushing existing rows 74-123 down to 76-125.
$ws.Rows("74:75").Insert()

# --- New row 74 ---
$ws.Cells.Item(74,1).Value  = 9
$ws.Cells.Item(74,2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(74,3).Value  = "Metropolitana"
$ws.Cells.Item(74,4).Value  = 44518
$ws.Cells.Item(74,5).Value  = 13
$ws.Cells.Item(74,6).Value  = "Fruta"
$ws.Cells.Item(74,7).Value  = 100101
$ws.Cells.Item(74,8).Value  = "Berries"
$ws.Cells.Item(74,9).Value  = 100101001
$ws.Cells.Item(74,10).Value = "Arándano (blue)"
$ws.Cells.Item(74,11).Value = "Sin especificar"
$ws.Cells.Item(74,12).Value = "Especial"
$ws.Cells.Item(74,13).Value = 330
$ws.Cells.Item(74,14).Value = 6000
$ws.Cells.Item(74,15).Value = 6000
$ws.Cells.Item(74,16).Value = 6000
$ws.Cells.Item(74,17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(74,18).Value = "Región Metropolitana"
$ws.Cells.Item(74,19).Value = 3000
$ws.Cells.Item(74,20).Value = 2

# --- New row 75 ---
$ws.Cells.Item(75,1).Value  = 9
$ws.Cells.Item(75,2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(75,3).Value  = "Metropolitana"
$ws.Cells.Item(75,4).Value  = 44518
$ws.Cells.Item(75,5).Value  = 13
$ws.Cells.Item(75,6).Value  = "Fruta"
$ws.Cells.Item(75,7).Value  = 100101
$ws.Cells.Item(75,8).Value  = "Berries"
$ws.Cells.Item(75,9).Value  = 100101001
$ws.Cells.Item(75,10).Value = "Arándano (blue)"
$ws.Cells.Item(75,11).Value = "Sin especificar"
$ws.Cells.Item(75,12).Value = "Primera"
$ws.Cells.Item(75,13).Value = 450
$ws.Cells.Item(75,14).Value = 5000
$ws.Cells.Item(75,15).Value = 5000
$ws.Cells.Item(75,16).Value = 5000
$ws.Cells.Item(75,17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(75,18).Value = "Región Metropolitana"
$ws.Cells.Item(75,19).Value = 2500
$ws.Cells.Item(75,20).Value = 2
